# Bug Report workbook update: add "Title" column and rework the old
# "Priority" header into "Environment details" on the Login sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# --- Insert a new column B for "Title" (shifts Description..Priority right) ---
$ws.Columns.Item(2).Insert()

$titleCell = $ws.Cells.Item(1, 2)
$titleCell.Value2 = "Title"
$titleCell.Font.Bold = $true
$titleCell.Font.Size = 11
$titleCell.Font.Name = "Calibri"
$titleCell.Interior.Color = 8242323
$titleCell.HorizontalAlignment = -4108

# --- The old "Priority" header is now column G; replace it with
#     "Environment details " and give it the same header styling ---
$envCell = $ws.Cells.Item(1, 7)
$envCell.Value2 = "Environment details "
$envCell.Font.Bold = $true
$envCell.Font.Size = 11
$envCell.Font.Name = "Calibri"
$envCell.Interior.Color = 8242323
$envCell.HorizontalAlignment = -4108

# The trailing space was typed/formatted separately at a larger, non-bold
# size (matches the rich-text run captured in the source workbook).
$trailingSpace = $envCell.Characters(20, 1)
$trailingSpace.Font.Size = 12
$trailingSpace.Font.Bold = $false
$trailingSpace.Font.Name = "Calibri"

# Widen the new columns a bit (matches authoring intent for the longer headers)
$ws.Columns.Item(7).ColumnWidth = 25.6
$ws.Columns.Item(8).ColumnWidth = 15.3

Write-Output "done"
